$d = $word.ActiveDocument

# Locate the paragraph that currently ends with the "55555555555555555555"
# run (and, in the original document, also carries the trailing
# _GoBack bookmark immediately after that run).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*55555555555555555555*") {
        $targetPara = $p
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing '55555555555555555555'"
}

# Remove the pre-existing _GoBack bookmark (it currently sits right after
# the "55555555555555555555" text, inside that same paragraph); it will be
# re-created at the correct spot - after the new "66666" run - below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapsed insertion point: right before the paragraph mark that ends the
# "55555555555555555555" paragraph.
$insertAt = $targetPara.Range.End - 1
$ins = $d.Range($insertAt, $insertAt)

# Insert a brand-new paragraph (with its own run "66666" and the _GoBack
# bookmark after it) via raw WordprocessingML so the exact run/paragraph
# formatting matches what Word itself would author.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>66666</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$ins.InsertXML($xmlFrag)
